$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "传记" (biography) entry row (sequence id 4) from the role/button
# list, which shifts all subsequent rows up by one.
$ws.Rows("6").Delete()

# Match the author's final cell selection.
$ws.Range("C7").Select() | Out-Null
